$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.845.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "'1.881.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.46%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'326.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").Value = "'0.4675"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.18%  "

$ws.Range("E8").Value = "  +2.69%  "

$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("D10").Value = "'0.9815"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.39%  "

$ws.Range("D11").Value = "'22.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.81%  "

$ws.Range("D12").Value = "'1.883.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.51%  "

$ws.Range("D13").Value = "'5.751"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("D14").Value = "'7.022"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.66%  "

$ws.Range("D15").Value = "'0.06973"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.73%  "

$ws.Range("E16").Value = "  +2.48%  "

$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").Value = "'0.00001009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.84%  "

$ws.Range("D19").Value = "'17.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").Value = "'1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("D21").Value = "'28.861.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.54%  "

$ws.Range("D22").Value = "'5.368"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("D23").Value = "'11.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.71%  "

$ws.Range("D24").Value = "'2.124"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.14%  "

$ws.Range("D25").Value = "'2.119.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.48%  "

$ws.Range("D26").Value = "'153.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.97%  "

$ws.Range("D27").Value = "'19.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.22%  "

$ws.Range("D28").Value = "'5.770"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").Value = "'2.008"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.29%  "

$ws.Range("D30").Value = "'120.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.24%  "

$ws.Range("D31").Value = "'0.09402"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.83%  "

$ws.Range("D32").Value = "'0.9425"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.11%  "

$ws.Range("D33").Value = "'5.321"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("E34").Value = "  +3.54%  "

$ws.Range("D35").Value = "'3.355"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.42%  "

$ws.Range("D36").Value = "'0.05924"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").Value = "'0.02123"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.67%  "

$ws.Range("D38").Value = "'1.151"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").Value = "'7.921"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.61%  "

$ws.Range("D40").Value = "'0.5730"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.12%  "

$ws.Range("D41").Value = "'0.1799"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.07%  "

$ws.Range("E42").Value = "  +1.33%  "

$ws.Range("D43").Value = "'0.07276"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.31%  "

$ws.Range("D44").Value = "'11.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("D45").Value = "'0.5347"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.25%  "

$ws.Range("D46").Value = "'1.151"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.34%  "

$ws.Range("D47").Value = "'2.124"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.30%  "

$ws.Range("D48").Value = "'1.853"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.13%  "

$ws.Range("D49").Value = "'114.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.99%  "

$ws.Range("D50").Value = "'2.370"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.43%  "

$ws.Range("D51").Value = "'1.007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.40%  "

